$d = $word.ActiveDocument

# 1) Sprint 1 user story points: Cognito Authentication (8 points) -> (11 points)
$d.Content.Find.Execute("ation (8 points)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ation (11 points)", 2) | Out-Null

# 2) Sprint 1 burndown total: Total: 16 Story Points -> Total: 19 Story Points
$d.Content.Find.Execute("Total: 16 Story Points", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Total: 19 Story Points", 2) | Out-Null

# 3) Sprint 2 burndown total: Total: 11 Story Points -> Total: 8 Story Points
$d.Content.Find.Execute("Total: 11 Story Points", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Total: 8 Story Points", 2) | Out-Null
